# Add season-record columns (Wins, Losses, Ties) after the existing
# last column (AC) of the roster/statistics table.
#
# - AD1/AE1/AF1 get the same header formatting as the existing header
#   row (bold, centered, bordered) by copying the format from AC1 and
#   then overwriting the text.
# - AD2:AF44 get the season record values (90 wins, 72 losses, 0 ties)
#   for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ------------------------------------------------------
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows --------------------------------------------------------
$wins = 90
$losses = 72
$ties = 0

for ($row = 2; $row -le 44; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins     # column AD
    $ws.Cells.Item($row, 31).Value = $losses   # column AE
    $ws.Cells.Item($row, 32).Value = $ties     # column AF
}
